{"js": "// Resume edit: change the \"Seeking role...\" objective sentence and the\n// internship end date to reflect a summer role at an asset management firm.\n//\n//   \"Seeking role in product/program management ...\"\n//     -> \"Seeking role in asset management firm ...\"\n//   \"... to help the team design and scale products and initiatives.\"\n//     -> \"... to help the team conduct analysis, build models, and drive\n//         profitable investments.\"\n//   \"June 2019 - Present (Extended) - San Jose\"\n//     -> \"June 2019 - March 2020 - San Jose\"\n\nasync function replaceOnce(searchText, replacement) {\n  const results = context.document.body.search(searchText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(replacement, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\nawait replaceOnce(\n  \"Seeking role in product/program management\",\n  \"Seeking role in asset management firm\"\n);\n\nawait replaceOnce(\n  \"design and scale products and initiatives\",\n  \"conduct analysis, build models, and drive profitable investments\"\n);\n\nawait replaceOnce(\"Present (Extended) \", \"March 2020 \");\n", "ps1": "# Resume edit: change the \"Seeking role...\" objective sentence and the\n# internship end date to reflect a summer role at an asset management firm.\n#\n#   \"Seeking role in product/program management ...\"\n#     -> \"Seeking role in asset management firm ...\"\n#   \"... to help the team design and scale products and initiatives.\"\n#     -> \"... to help the team conduct analysis, build models, and drive\n#         profitable investments.\"\n#   \"June 2019 - Present (Extended) - San Jose\"\n#     -> \"June 2019 - March 2020 - San Jose\"\n\n$d = $word.ActiveDocument\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nfunction Replace-Text($searchText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($searchText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll) | Out-Null\n}\n\nReplace-Text \"Seeking role in product/program management\" \"Seeking role in asset management firm\"\nReplace-Text \"design and scale products and initiatives\" \"conduct analysis, build models, and drive profitable investments\"\nReplace-Text \"Present (Extended) \" \"March 2020 \"\n"}
